$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 22.11125233333334
$ws.Range("H2").Value = 66.33375700000001
$ws.Range("I2").Value = 0.08763778737242772
$ws.Range("J2").Value = 0.08763778737242772
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.024701
$ws.Range("N2").Value = 9.074103000000001
$ws.Range("O2").Value = 0.1596375877334842
$ws.Range("P2").Value = 0.1596375877334843
$ws.Range("Q2").Value = 66.87992704388569
$ws.Range("R2").Value = 601.9193433949711
$ws.Range("S2").Value = 0.01399028497043437
$ws.Range("T2").Value = 0.01399028497043437
$ws.Range("G3").Value = 22.11125233333334
$ws.Range("H3").Value = 66.33375700000001
$ws.Range("I3").Value = 0.08763778737242772
$ws.Range("J3").Value = 0.08763778737242772
$ws.Range("O3").Value = 0.6072559333217162
$ws.Range("P3").Value = 0.6072559333217163
$ws.Range("Q3").Value = 254.4089590311718
$ws.Range("R3").Value = 2289.680631280546
$ws.Range("S3").Value = 0.05321856636509371
$ws.Range("T3").Value = 0.05321856636509372
$ws.Range("G4").Value = 22.11125233333334
$ws.Range("H4").Value = 66.33375700000001
$ws.Range("I4").Value = 0.08763778737242772
$ws.Range("J4").Value = 0.08763778737242772
$ws.Range("M4").Value = 4.368554666666666
$ws.Range("N4").Value = 13.105664
$ws.Range("O4").Value = 0.2305634602787257
$ws.Range("P4").Value = 0.2305634602787257
$ws.Range("Q4").Value = 96.59421456662756
$ws.Range("R4").Value = 869.347931099648
$ws.Range("S4").Value = 0.02020607150775815
$ws.Range("T4").Value = 0.02020607150775815
$ws.Range("G5").Value = 22.11125233333334
$ws.Range("H5").Value = 66.33375700000001
$ws.Range("I5").Value = 0.08763778737242772
$ws.Range("J5").Value = 0.08763778737242772
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.04818333333333333
$ws.Range("N5").Value = 0.14455
$ws.Range("O5").Value = 0.002543018666073676
$ws.Range("P5").Value = 0.002543018666073677
$ws.Range("Q5").Value = 1.065393841594445
$ws.Range("R5").Value = 9.588544574350001
$ws.Range("S5").Value = 0.0002228645291414796
$ws.Range("T5").Value = 0.0002228645291414797
$ws.Range("I6").Value = 0.8254813202458152
$ws.Range("J6").Value = 0.8254813202458152
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.024701
$ws.Range("N6").Value = 9.074103000000001
$ws.Range("O6").Value = 0.1596375877334842
$ws.Range("P6").Value = 0.1596375877334843
$ws.Range("Q6").Value = 629.9580595242177
$ws.Range("R6").Value = 5669.622535717959
$ws.Range("S6").Value = 0.1317778466830937
$ws.Range("T6").Value = 0.1317778466830938
$ws.Range("I7").Value = 0.8254813202458152
$ws.Range("J7").Value = 0.8254813202458152
$ws.Range("O7").Value = 0.6072559333217162
$ws.Range("P7").Value = 0.6072559333217163
$ws.Range("S7").Value = 0.501278429565515
$ws.Range("T7").Value = 0.5012784295655152
$ws.Range("I8").Value = 0.8254813202458152
$ws.Range("J8").Value = 0.8254813202458152
$ws.Range("M8").Value = 4.368554666666666
$ws.Range("N8").Value = 13.105664
$ws.Range("O8").Value = 0.2305634602787257
$ws.Range("P8").Value = 0.2305634602787257
$ws.Range("Q8").Value = 909.8440542515767
$ws.Range("R8").Value = 8188.596488264191
$ws.Range("S8").Value = 0.1903258295913261
$ws.Range("T8").Value = 0.1903258295913261
$ws.Range("I9").Value = 0.8254813202458152
$ws.Range("J9").Value = 0.8254813202458152
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.04818333333333333
$ws.Range("N9").Value = 0.14455
$ws.Range("O9").Value = 0.002543018666073676
$ws.Range("P9").Value = 0.002543018666073677
$ws.Range("Q9").Value = 10.03519989846111
$ws.Range("R9").Value = 90.31679908615
$ws.Range("S9").Value = 0.00209921440588025
$ws.Range("T9").Value = 0.002099214405880251
$ws.Range("G10").Value = 12.43397833333333
$ws.Range("H10").Value = 37.301935
$ws.Range("I10").Value = 0.04928198244688778
$ws.Range("J10").Value = 0.04928198244688778
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.024701
$ws.Range("N10").Value = 9.074103000000001
$ws.Range("O10").Value = 0.1596375877334842
$ws.Range("P10").Value = 0.1596375877334843
$ws.Range("Q10").Value = 37.60906669881167
$ws.Range("R10").Value = 338.481600289305
$ws.Range("S10").Value = 0.007867256796545078
$ws.Range("T10").Value = 0.007867256796545079
$ws.Range("G11").Value = 12.43397833333333
$ws.Range("H11").Value = 37.301935
$ws.Range("I11").Value = 0.04928198244688778
$ws.Range("J11").Value = 0.04928198244688778
$ws.Range("O11").Value = 0.6072559333217162
$ws.Range("P11").Value = 0.6072559333217163
$ws.Range("Q11").Value = 143.0636056570478
$ws.Range("R11").Value = 1287.57245091343
$ws.Range("S11").Value = 0.02992677624672927
$ws.Range("T11").Value = 0.02992677624672928
$ws.Range("G12").Value = 12.43397833333333
$ws.Range("H12").Value = 37.301935
$ws.Range("I12").Value = 0.04928198244688778
$ws.Range("J12").Value = 0.04928198244688778
$ws.Range("M12").Value = 4.368554666666666
$ws.Range("N12").Value = 13.105664
$ws.Range("O12").Value = 0.2305634602787257
$ws.Range("P12").Value = 0.2305634602787257
$ws.Range("Q12").Value = 54.31851407331555
$ws.Range("R12").Value = 488.86662665984
$ws.Range("S12").Value = 0.01136262440234987
$ws.Range("T12").Value = 0.01136262440234987
$ws.Range("G13").Value = 12.43397833333333
$ws.Range("H13").Value = 37.301935
$ws.Range("I13").Value = 0.04928198244688778
$ws.Range("J13").Value = 0.04928198244688778
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.04818333333333333
$ws.Range("N13").Value = 0.14455
$ws.Range("O13").Value = 0.002543018666073676
$ws.Range("P13").Value = 0.002543018666073677
$ws.Range("Q13").Value = 0.5991105226944445
$ws.Range("R13").Value = 5.39199470425
$ws.Range("S13").Value = 0.0001253250012635509
$ws.Range("T13").Value = 0.0001253250012635509
$ws.Range("G14").Value = 9.486307333333334
$ws.Range("H14").Value = 28.458922
$ws.Range("I14").Value = 0.03759890993486929
$ws.Range("J14").Value = 0.03759890993486929
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 3.024701
$ws.Range("N14").Value = 9.074103000000001
$ws.Range("O14").Value = 0.1596375877334842
$ws.Range("P14").Value = 0.1596375877334843
$ws.Range("Q14").Value = 28.69324327744067
$ws.Range("R14").Value = 258.239189496966
$ws.Range("S14").Value = 0.006002199283411069
$ws.Range("T14").Value = 0.00600219928341107
$ws.Range("G15").Value = 9.486307333333334
$ws.Range("H15").Value = 28.458922
$ws.Range("I15").Value = 0.03759890993486929
$ws.Range("J15").Value = 0.03759890993486929
$ws.Range("O15").Value = 0.6072559333217162
$ws.Range("P15").Value = 0.6072559333217163
$ws.Range("Q15").Value = 109.1481177701018
$ws.Range("R15").Value = 982.3330599309161
$ws.Range("S15").Value = 0.0228321611443782
$ws.Range("T15").Value = 0.0228321611443782
$ws.Range("G16").Value = 9.486307333333334
$ws.Range("H16").Value = 28.458922
$ws.Range("I16").Value = 0.03759890993486929
$ws.Range("J16").Value = 0.03759890993486929
$ws.Range("M16").Value = 4.368554666666666
$ws.Range("N16").Value = 13.105664
$ws.Range("O16").Value = 0.2305634602787257
$ws.Range("P16").Value = 0.2305634602787257
$ws.Range("Q16").Value = 41.44145217046756
$ws.Range("R16").Value = 372.973069534208
$ws.Range("S16").Value = 0.008668934777291621
$ws.Range("T16").Value = 0.008668934777291623
$ws.Range("G17").Value = 9.486307333333334
$ws.Range("H17").Value = 28.458922
$ws.Range("I17").Value = 0.03759890993486929
$ws.Range("J17").Value = 0.03759890993486929
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.04818333333333333
$ws.Range("N17").Value = 0.14455
$ws.Range("O17").Value = 0.002543018666073676
$ws.Range("P17").Value = 0.002543018666073677
$ws.Range("Q17").Value = 0.4570819083444445
$ws.Range("R17").Value = 4.113737175100001
$ws.Range("S17").Value = 0.0000956147297883956
$ws.Range("T17").Value = 0.00009561472978839563
